$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.247.12"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "2.836.15"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'361.73"
$ws.Range("E5").Value = "  +3.07%  "
$ws.Range("D6").Value = "'111.91"
$ws.Range("E6").Value = "  -3.72%  "
$ws.Range("D7").Value = "'0.568"
$ws.Range("E7").Value = "  +3.30%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D9").Value = "'0.603"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("D10").Value = "'40.89"
$ws.Range("E10").Value = "  -3.99%  "
$ws.Range("D11").Value = "'0.0865"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").Value = "'19.97"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").Value = "3.287.45"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "2.844.95"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").Value = "'0.929"
$ws.Range("E17").Value = "  +3.66%  "
$ws.Range("D18").Value = "52.135.38"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "'7.49"
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").Value = "'13.39"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "0.0₃0999"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("D23").Value = "'272.16"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").Value = "'70.36"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'2.81"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("D26").Value = "'26.91"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "'10.31"
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("D29").Value = "'2.24"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("E30").Value = "  +1.63%  "
$ws.Range("E31").Value = "  +4.60%  "
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").Value = "'52.40"
$ws.Range("E32").Value = "  +3.92%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "'35.06"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("D34").Value = "'5.88"
$ws.Range("E34").Value = "  +1.07%  "
$ws.Range("D35").Value = "'5.57"
$ws.Range("E35").Value = "  +12.22%  "
$ws.Range("D36").Value = "'0.0853"
$ws.Range("E36").Value = "  +1.80%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "'3.27"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").Value = "'2.04"
$ws.Range("E39").Value = "  -3.46%  "
$ws.Range("D40").Value = "'18.37"
$ws.Range("E40").Value = "  -2.30%  "
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("D42").Value = "'2.54"
$ws.Range("E42").Value = "  -2.87%  "
$ws.Range("D43").Value = "'125.27"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'2.27"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'22.55"
$ws.Range("E45").Value = "  -5.69%  "
$ws.Range("D46").Value = "2.082.88"
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("D47").Value = "'3.32"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("D48").Value = "'2.32"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("D49").Value = "'5.87"
$ws.Range("E49").Value = "  +5.15%  "
$ws.Range("D50").Value = "'0.968"
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("D51").Value = "'9.20"
$ws.Range("E51").Value = "  +2.20%  "
